# Apply updated metrics to rows 2-26 (all data rows share the same new values)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B" = 0.9999883732771242
    "C" = 0.9990763349190175
    "D" = 0.9999960854524444
    "E" = 0.9999999999998683
    "F" = 0.9999980425969037
    "G" = 0.00001085303545036088
    "H" = 0.000862200808886027
    "I" = 0.000003020773121276768
    "J" = 0.0000000000001009460136362478
    "K" = 0.000001510386611111391
    "L" = 0.0001805656912053873
    "M" = 0.00329439454989242
    "N" = 0.9999069862169936
    "O" = 0.003434643746218193
    "P" = 64.86213150418904
    "Q" = 90.45852382642124
}

for ($row = 2; $row -le 26; $row++) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
